# Workbook "Hortaliza, Vega Modelo de Temuco - Ajo": weekly fruit/vegetable price update.
# Two new daily price observations (rows) are inserted at the top of the
# existing "Ajo" (garlic) data block, pushing the rest of the data set
# down by two rows (old row 1239 becomes 1241, ..., old row 1305 becomes 1307).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 1239; this shifts all
# subsequent rows (1239..1305) down to (1241..1307), preserving their
# values and formatting intact, and extends the sheet dimension to A1:R1307.
$ws.Rows("1239:1240").Insert()

# Fill in the first new row (1239) with the new observation.
$ws.Cells.Item(1239, 1).Value = 10
$ws.Cells.Item(1239, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1239, 3).Value = "La Araucanía"
$ws.Cells.Item(1239, 4).Value = 45267
$ws.Cells.Item(1239, 5).Value = 9
$ws.Cells.Item(1239, 6).Value = 100112003
$ws.Cells.Item(1239, 7).Value = "Ajo"
$ws.Cells.Item(1239, 8).Value = "Chino"
$ws.Cells.Item(1239, 9).Value = "Primera"
$ws.Cells.Item(1239, 10).Value = 650
$ws.Cells.Item(1239, 11).Value = 24000
$ws.Cells.Item(1239, 12).Value = 25000
$ws.Cells.Item(1239, 13).Value = 24615
$ws.Cells.Item(1239, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(1239, 15).Value = "China"
$ws.Cells.Item(1239, 16).Value = 2462
$ws.Cells.Item(1239, 17).Value = 10
$ws.Cells.Item(1239, 18).Value = "Hortaliza"

# Fill in the second new row (1240) with the new observation.
$ws.Cells.Item(1240, 1).Value = 10
$ws.Cells.Item(1240, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1240, 3).Value = "La Araucanía"
$ws.Cells.Item(1240, 4).Value = 45267
$ws.Cells.Item(1240, 5).Value = 9
$ws.Cells.Item(1240, 6).Value = 100112003
$ws.Cells.Item(1240, 7).Value = "Ajo"
$ws.Cells.Item(1240, 8).Value = "Chino"
$ws.Cells.Item(1240, 9).Value = "Primera"
$ws.Cells.Item(1240, 10).Value = 400
$ws.Cells.Item(1240, 11).Value = 26000
$ws.Cells.Item(1240, 12).Value = 27000
$ws.Cells.Item(1240, 13).Value = 26250
$ws.Cells.Item(1240, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(1240, 15).Value = "China"
$ws.Cells.Item(1240, 16).Value = 2625
$ws.Cells.Item(1240, 17).Value = 10
$ws.Cells.Item(1240, 18).Value = "Hortaliza"
